$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CL")

# Row 6: Change in inventories
$ws.Range("B6").Value = -320000000.0
$ws.Range("C6").Value = -251000000.0
$ws.Range("D6").Value = -222000000.0
$ws.Range("E6").Value = -190000000.0
$ws.Range("F6").Value = -16000000.0
$ws.Range("G6").Value = -77000000.0

# Row 7: Change in payables and accrued liability
$ws.Range("B7").Value = 225000000.0
$ws.Range("C7").Value = 520000000.0
$ws.Range("D7").Value = 556000000.0
$ws.Range("E7").Value = 397000000.0
$ws.Range("F7").Value = 212000000.0
$ws.Range("G7").Value = 36000000.0

# Row 17: Dividends Paid (Total)
$ws.Range("B17").Value = -1505000000.0

# Row 25: Dividends Paid (Common)
$ws.Range("B25").Value = -1657000000.0

# Row 29: Capital Stock Change
$ws.Range("B29").Value = -1021000000.0
